$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Copy formatting from column Q to column R for rows 3 through 37,
# then set the new (2021) values. Row 34 stays format-only (no value),
# matching its blank Q34 sibling.
for ($r = 3; $r -le 37; $r++) {
    $ws.Range("Q$r").Copy() | Out-Null
    $ws.Range("R$r").PasteSpecial($xlPasteFormats) | Out-Null
}
$excel.CutCopyMode = $false

# Set the 2021 data values that were added in column R.
$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 0.12641839647678207
$ws.Range("R5").Value = 0.14922981985616976
$ws.Range("R6").Value = 0.10326895933792253
$ws.Range("R7").Value = [double]"3.433011112114915E-2"
$ws.Range("R8").Value = [double]"3.6820478077087354E-2"
$ws.Range("R9").Value = [double]"3.1930519190242035E-2"
$ws.Range("R10").Value = [double]"8.7302929367211068E-2"
$ws.Range("R11").Value = 0.10296328329317765
$ws.Range("R12").Value = [double]"7.1859056271889668E-2"
$ws.Range("R13").Value = 0.10716050460690947
$ws.Range("R14").Value = [double]"7.9035451351703812E-2"
$ws.Range("R15").Value = 0.13553052227085377
$ws.Range("R16").Value = [double]"6.479643687803946E-2"
$ws.Range("R17").Value = [double]"7.643825526207898E-2"
$ws.Range("R18").Value = [double]"5.3576570965516782E-2"
$ws.Range("R19").Value = [double]"5.4163459619715498E-2"
$ws.Range("R20").Value = [double]"6.4872252119520635E-2"
$ws.Range("R21").Value = [double]"4.3693418784505472E-2"
$ws.Range("R22").Value = [double]"5.1373884452794741E-2"
$ws.Range("R23").Value = [double]"2.9662368095156877E-2"
$ws.Range("R24").Value = [double]"7.2642215296997686E-2"
$ws.Range("R25").Value = 0.13772601093442507
$ws.Range("R26").Value = 0.15668565643254884
$ws.Range("R27").Value = 0.11816042869432726
$ws.Range("R28").Value = 0.33417383115107696
$ws.Range("R29").Value = 0.41139191068108794
$ws.Range("R30").Value = 0.24697746624641295
$ws.Range("R31").Value = 0.16773611144997194
$ws.Range("R32").Value = 0.1959922553363346
$ws.Range("R33").Value = 0.13791201213625709
$ws.Range("R35").Value = 0
$ws.Range("R36").Value = 0.1
$ws.Range("R37").Value = 0.2

# Update the active selection (sheetView) to C1.
$ws.Range("C1").Select()

Write-Output "R column (2021) populated."
